$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1961904761904762
$ws.Range("C2").Value = 0.5428571428571428
$ws.Range("J2").Value = 0.009523809523809525
$ws.Range("P2").Value = 0.1523809523809524
$ws.Range("S2").Value = 0.09904761904761905
$ws.Range("B3").Value = 0.006666666666666667
$ws.Range("C3").Value = 0.03666666666666667
$ws.Range("J3").Value = 0.02333333333333333
$ws.Range("P3").Value = 0.7433333333333333
$ws.Range("S3").Value = 0.19
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.7285714285714285
$ws.Range("S4").Value = 0.2
$ws.Range("P5").Value = 0.8333333333333334
$ws.Range("S5").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.06310679611650485
$ws.Range("D6").Value = 0.007281553398058253
$ws.Range("F6").Value = 0.03640776699029126
$ws.Range("J6").Value = 0.3009708737864077
$ws.Range("O6").Value = 0.01456310679611651
$ws.Range("Q6").Value = 0.1771844660194175
$ws.Range("R6").Value = 0.06067961165048544
$ws.Range("S6").Value = 0.3398058252427185
$ws.Range("B7").Value = 0.1082474226804124
$ws.Range("D7").Value = 0.02577319587628866
$ws.Range("F7").Value = 0.05670103092783505
$ws.Range("J7").Value = 0.1314432989690722
$ws.Range("O7").Value = 0.02061855670103093
$ws.Range("Q7").Value = 0.1907216494845361
$ws.Range("R7").Value = 0.06701030927835051
$ws.Range("S7").Value = 0.3994845360824743
$ws.Range("B8").Value = 0.08695652173913043
$ws.Range("D8").Value = 0.01484623541887593
$ws.Range("E8").Value = 0.002120890774125133
$ws.Range("F8").Value = 0.0784729586426299
$ws.Range("J8").Value = 0.1039236479321315
$ws.Range("O8").Value = 0.01696712619300106
$ws.Range("Q8").Value = 0.1919406150583245
$ws.Range("R8").Value = 0.0784729586426299
$ws.Range("S8").Value = 0.4262990455991517
$ws.Range("B9").Value = 0.1187335092348285
$ws.Range("D9").Value = 0.0079155672823219
$ws.Range("E9").Value = 0.002638522427440633
$ws.Range("F9").Value = 0.0395778364116095
$ws.Range("J9").Value = 0.1187335092348285
$ws.Range("O9").Value = 0.02110817941952507
$ws.Range("Q9").Value = 0.1767810026385224
$ws.Range("R9").Value = 0.06596306068601583
$ws.Range("S9").Value = 0.4485488126649076
$ws.Range("B10").Value = 0.09904153354632587
$ws.Range("D10").Value = 0.01962574167047011
$ws.Range("E10").Value = 0.001825650387950707
$ws.Range("F10").Value = 0.07074395253308992
$ws.Range("J10").Value = 0.1077133728890917
$ws.Range("O10").Value = 0.01232314011866728
$ws.Range("Q10").Value = 0.2117754450022821
$ws.Range("R10").Value = 0.07941579187585578
$ws.Range("S10").Value = 0.3975353719762665
$ws.Range("G11").Value = 0.1341222879684418
$ws.Range("J11").Value = 0.0650887573964497
$ws.Range("K11").Value = 0.1577909270216963
$ws.Range("L11").Value = 0.631163708086785
$ws.Range("S11").Value = 0.01183431952662722
$ws.Range("G12").Value = 0.7586206896551724
$ws.Range("J12").Value = 0.1810344827586207
$ws.Range("K12").Value = 0.005747126436781609
$ws.Range("L12").Value = 0.02298850574712644
$ws.Range("S12").Value = 0.03160919540229885
$ws.Range("G13").Value = 0.7291666666666666
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.02083333333333333
$ws.Range("G14").Value = 0.5
$ws.Range("S14").Value = 0.5
$ws.Range("F15").Value = 0.01699029126213592
$ws.Range("H15").Value = 0.1723300970873786
$ws.Range("I15").Value = 0.08495145631067962
$ws.Range("J15").Value = 0.3398058252427185
$ws.Range("K15").Value = 0.0703883495145631
$ws.Range("M15").Value = 0.01941747572815534
$ws.Range("O15").Value = 0.0558252427184466
$ws.Range("S15").Value = 0.2402912621359223
$ws.Range("F16").Value = 0.01704545454545454
$ws.Range("H16").Value = 0.21875
$ws.Range("I16").Value = 0.09659090909090909
$ws.Range("J16").Value = 0.34375
$ws.Range("K16").Value = 0.1079545454545455
$ws.Range("M16").Value = 0.04261363636363636
$ws.Range("O16").Value = 0.07670454545454546
$ws.Range("S16").Value = 0.09659090909090909
$ws.Range("F17").Value = 0.02100350058343057
$ws.Range("H17").Value = 0.2042007001166861
$ws.Range("I17").Value = 0.07934655775962661
$ws.Range("J17").Value = 0.3955659276546091
$ws.Range("K17").Value = 0.102683780630105
$ws.Range("M17").Value = 0.01750291715285881
$ws.Range("O17").Value = 0.06767794632438739
$ws.Range("S17").Value = 0.1120186697782964
$ws.Range("F18").Value = 0.009375
$ws.Range("H18").Value = 0.221875
$ws.Range("I18").Value = 0.05
$ws.Range("J18").Value = 0.390625
$ws.Range("K18").Value = 0.1375
$ws.Range("M18").Value = 0.025
$ws.Range("N18").Value = 0.003125
$ws.Range("O18").Value = 0.06875000000000001
$ws.Range("S18").Value = 0.09375
$ws.Range("F19").Value = 0.02259164535379369
$ws.Range("H19").Value = 0.23231031543052
$ws.Range("I19").Value = 0.09505541346973571
$ws.Range("J19").Value = 0.3554987212276215
$ws.Range("K19").Value = 0.09335038363171355
$ws.Range("M19").Value = 0.02216538789428815
$ws.Range("N19").Value = 0.0004262574595055413
$ws.Range("O19").Value = 0.07289002557544758
$ws.Range("S19").Value = 0.1057118499573743
